$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 25
$ws.Range("K2").Value = 41
$ws.Range("E3").Value = 46
$ws.Range("G3").Value = 41
$ws.Range("H3").Value = 30
$ws.Range("L3").Value = 71
$ws.Range("I6").Value = 10
$ws.Range("C9").Value = 149
$ws.Range("D9").Value = 147
$ws.Range("E9").Value = 152
$ws.Range("F9").Value = 175
$ws.Range("I9").Value = 177
$ws.Range("J9").Value = 132
$ws.Range("K9").Value = 136
$ws.Range("L9").Value = 156
$ws.Range("B10").Value = 330
$ws.Range("C10").Value = 419
$ws.Range("D10").Value = 555
$ws.Range("E10").Value = 614
$ws.Range("F10").Value = 736
$ws.Range("G10").Value = 480
$ws.Range("I10").Value = 265
$ws.Range("J10").Value = 232
$ws.Range("L10").Value = 222
$ws.Range("B11").Value = 503
$ws.Range("C11").Value = 619
$ws.Range("D11").Value = 781
$ws.Range("E11").Value = 838
$ws.Range("F11").Value = 980
$ws.Range("G11").Value = 735
$ws.Range("H11").Value = 357
$ws.Range("I11").Value = 540
$ws.Range("J11").Value = 483
$ws.Range("K11").Value = 479
$ws.Range("L11").Value = 505

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F5").Value = 11
$ws.Range("D7").Value = 9
$ws.Range("H7").Value = 5
$ws.Range("F8").Value = 41
$ws.Range("L8").Value = 18
$ws.Range("D13").Value = 4
$ws.Range("J15").Value = 3
$ws.Range("K18").Value = 12
$ws.Range("L18").Value = 14
$ws.Range("C19").Value = 5
$ws.Range("C27").Value = 44
$ws.Range("D27").Value = 32
$ws.Range("E27").Value = 29
$ws.Range("D28").Value = 12
$ws.Range("D31").Value = 34
$ws.Range("L31").Value = 27
$ws.Range("C35").Value = 26
$ws.Range("K35").Value = 25
$ws.Range("G40").Value = 11
$ws.Range("B44").Value = 6
$ws.Range("E44").Value = 7
$ws.Range("C46").Value = 21
$ws.Range("D46").Value = 16
$ws.Range("G46").Value = 26
$ws.Range("I48").Value = 4
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = 26
$ws.Range("L49").Value = 16
$ws.Range("C51").Value = 8
$ws.Range("L51").Value = 9
$ws.Range("C52").Value = 84
$ws.Range("D52").Value = 199
$ws.Range("E52").Value = 216
$ws.Range("F52").Value = 245
$ws.Range("I52").Value = 116
$ws.Range("K52").Value = 68
$ws.Range("C60").Value = 10
$ws.Range("L60").Value = 2
$ws.Range("E62").Value = 2
$ws.Range("B64").Value = 13
$ws.Range("K69").Value = 12
$ws.Range("E73").Value = 20
$ws.Range("C75").Value = 24
$ws.Range("I75").Value = 12
$ws.Range("J75").Value = 17
$ws.Range("L75").Value = 24
$ws.Range("I76").Value = 25
$ws.Range("E79").Value = 6
$ws.Range("G79").Value = 3
$ws.Range("C81").Value = 13
$ws.Range("D86").Value = 4
$ws.Range("F90").Value = 21
$ws.Range("E93").Value = 33
$ws.Range("B97").Value = 503
$ws.Range("C97").Value = 619
$ws.Range("D97").Value = 781
$ws.Range("E97").Value = 838
$ws.Range("F97").Value = 980
$ws.Range("G97").Value = 735
$ws.Range("H97").Value = 357
$ws.Range("I97").Value = 540
$ws.Range("J97").Value = 483
$ws.Range("K97").Value = 479
$ws.Range("L97").Value = 505

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("C7").Value = 20
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 13
$ws.Range("L7").Value = 15
$ws.Range("C8").Value = 24
$ws.Range("I8").Value = 12
$ws.Range("J8").Value = 17
$ws.Range("L8").Value = 24

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I7").Value = 8
$ws.Range("I9").Value = 25

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("H3").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 9
$ws.Range("H7").Value = 5

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F5").Value = 10
$ws.Range("L6").Value = 7
$ws.Range("F7").Value = 41
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 4

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 5

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L6").Value = 13
$ws.Range("D7").Value = 19
$ws.Range("D8").Value = 34
$ws.Range("L8").Value = 27

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 5
$ws.Range("K7").Value = 12
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 3
$ws.Range("C7").Value = 8
$ws.Range("C9").Value = 26
$ws.Range("K9").Value = 25

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I5").Value = 2
$ws.Range("K8").Value = 18
$ws.Range("C9").Value = 65
$ws.Range("D9").Value = 166
$ws.Range("E9").Value = 179
$ws.Range("F9").Value = 220
$ws.Range("C10").Value = 84
$ws.Range("D10").Value = 199
$ws.Range("E10").Value = 216
$ws.Range("F10").Value = 245
$ws.Range("I10").Value = 116
$ws.Range("K10").Value = 68

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F6").Value = 8
$ws.Range("F7").Value = 11

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 4
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 11
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 26
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("B5").Value = 9
$ws.Range("B6").Value = 13

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("G6").Value = 9
$ws.Range("G7").Value = 11

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("E4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("E6").Value = 6
$ws.Range("G6").Value = 3

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 5
$ws.Range("C5").Value = 31
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 32
$ws.Range("E6").Value = 29

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("G3").Value = 2
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 13
$ws.Range("C8").Value = 21
$ws.Range("D8").Value = 16
$ws.Range("G8").Value = 26

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("B5").Value = 6
$ws.Range("E5").Value = 6
$ws.Range("B6").Value = 6
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 12

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 13

$ws = $wb.Worksheets.Item("River North")
$ws.Range("E5").Value = 20
$ws.Range("E6").Value = 20

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("F6").Value = 18
$ws.Range("F7").Value = 21

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("C6").Value = 6
$ws.Range("L6").Value = 5
$ws.Range("C7").Value = 8
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 4

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("E5").Value = 31
$ws.Range("E6").Value = 33

$ws = $wb.Worksheets.Item("New City")
$ws.Range("E3").Value = 2
$ws.Range("E6").Value = 2

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("G3").Value = 2
$ws.Range("G5").Value = 3
